$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look like plain numbers
# so Excel does not auto-convert them from text to numeric type.
$textCells = @("D5","D6","D7","D9","D10","D12","D13","D14","D16","D21","D22","D23","D24","D27","D30","D31","D32","D34","D36","D37","D39","D41","D42","D43","D45","D46","D47","D48","D49","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply cell value updates
$ws.Range("D2").Value = "43.289.07"
$ws.Range("E2").Value = "  +1.04%  "
$ws.Range("D3").Value = "2.276.03"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "114.18"
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("D6").Value = "302.80"
$ws.Range("E6").Value = "  +6.33%  "
$ws.Range("D7").Value = "0.636"
$ws.Range("E7").Value = "  +1.31%  "
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("D9").Value = "0.616"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").Value = "44.75"
$ws.Range("E10").Value = "  -4.15%  "
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "55.24"
$ws.Range("E12").Value = "  +1.65%  "
$ws.Range("D13").Value = "8.89"
$ws.Range("E13").Value = "  -3.57%  "
$ws.Range("D14").Value = "1.05"
$ws.Range("E14").Value = "  +19.06%  "
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("D16").Value = "15.44"
$ws.Range("D17").Value = "2.620.51"
$ws.Range("E17").Value = "  +1.80%  "
$ws.Range("D18").Value = "2.278.64"
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("D19").Value = "43.208.77"
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").Value = "7.24"
$ws.Range("E21").Value = "  +5.31%  "
$ws.Range("D22").Value = "75.30"
$ws.Range("E22").Value = "  +4.12%  "
$ws.Range("D23").Value = "3.55"
$ws.Range("E23").Value = "  +11.25%  "
$ws.Range("D24").Value = "256.66"
$ws.Range("E24").Value = "  +10.58%  "
$ws.Range("E25").Value = "  +3.84%  "
$ws.Range("E26").Value = "  -1.81%  "
$ws.Range("D27").Value = "11.72"
$ws.Range("E27").Value = "  -2.64%  "
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "38.31"
$ws.Range("E30").Value = "  -5.32%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "175.36"
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "22.26"
$ws.Range("E32").Value = "  +5.16%  "
$ws.Range("E33").Value = "  -3.31%  "
$ws.Range("D34").Value = "0.0899"
$ws.Range("E34").Value = "  -0.36%  "
$ws.Range("E35").Value = "  +2.71%  "
$ws.Range("D36").Value = "5.08"
$ws.Range("E36").Value = "  +9.29%  "
$ws.Range("D37").Value = "4.29"
$ws.Range("E37").Value = "  -6.81%  "
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("D39").Value = "0.0378"
$ws.Range("E39").Value = "  +1.58%  "
$ws.Range("E40").Value = "  -0.81%  "
$ws.Range("D41").Value = "2.48"
$ws.Range("E41").Value = "  -5.09%  "
$ws.Range("D42").Value = "72.44"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").Value = "0.233"
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "12.71"
$ws.Range("E45").Value = "  -5.75%  "
$ws.Range("D46").Value = "1.36"
$ws.Range("E46").Value = "  +1.70%  "
$ws.Range("D47").Value = "5.64"
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("D48").Value = "107.44"
$ws.Range("E48").Value = "  +6.32%  "
$ws.Range("D49").Value = "1.30"
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("E50").Value = "  +2.76%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "1.67"
$ws.Range("E51").Value = "  +9.94%  "

# Restore default General number format on the forced-text cells
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "General"
}
